$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TFEC")

# Duplicate the formatting of row 6 into row 7 (shift-insert preserves the
# underlying style indices used throughout column A-C, including the
# column C style that is otherwise indistinguishable from "no style").
$ws.Range("A6:D6").Copy()
$ws.Range("A7:D7").Insert(-4121)
$excel.CutCopyMode = $false

# The insert-shift operation does not fully preserve the border formatting
# for columns A and B, so re-apply their exact formatting (reusing the
# existing styles) via a targeted format-only paste.
$ws.Range("A6").Copy()
$ws.Range("A7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("B6").Copy()
$ws.Range("B7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the new "Solar cookers" row values.
$ws.Range("A7").Value = "Solar cookers"
$ws.Range("B7").Value = "RES_CWH_SOLAR"
$ws.Range("C7").Value = "Solar stoves"
$ws.Range("D7").Value = "Solar"

$ws.Range("D8").Select()
